$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46070

# Row 3
$ws.Range("C3").Value = 46070

# Row 4
$ws.Range("C4").Value = 46070

# Row 5
$ws.Range("A5").Value = 'A 393-2025'
$ws.Range("B5").Value = 45663
$ws.Range("C5").Value = 46070
$ws.Range("G5").Value = 5.1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("R5").Value = 'Spillkråka'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/artfynd/A 393-2025 artfynd.xlsx", "A 393-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/kartor/A 393-2025 karta.png", "A 393-2025")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomål/A 393-2025 FSC-klagomål.docx", "A 393-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomålsmail/A 393-2025 FSC-klagomål mail.docx", "A 393-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsyn/A 393-2025 tillsynsbegäran.docx", "A 393-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsynsmail/A 393-2025 tillsynsbegäran mail.docx", "A 393-2025")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/fåglar/A 393-2025 prioriterade fågelarter.docx", "A 393-2025")'

# Row 6
$ws.Range("A6").Value = 'A 49789-2023'
$ws.Range("B6").Value = 45212
$ws.Range("C6").Value = 46070
$ws.Range("G6").Value = 3.8
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("R6").Value = 'Trubbfjädermossa'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/artfynd/A 49789-2023 artfynd.xlsx", "A 49789-2023")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/kartor/A 49789-2023 karta.png", "A 49789-2023")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomål/A 49789-2023 FSC-klagomål.docx", "A 49789-2023")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomålsmail/A 49789-2023 FSC-klagomål mail.docx", "A 49789-2023")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsyn/A 49789-2023 tillsynsbegäran.docx", "A 49789-2023")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsynsmail/A 49789-2023 tillsynsbegäran mail.docx", "A 49789-2023")'
$ws.Range("Z6").ClearContents()

# Row 7
$ws.Range("C7").Value = 46070

# Row 8
$ws.Range("C8").Value = 46070

# Row 9
$ws.Range("C9").Value = 46070

# Row 10
$ws.Range("C10").Value = 46070

# Row 11
$ws.Range("C11").Value = 46070

# Row 12
$ws.Range("C12").Value = 46070

# Row 13
$ws.Range("C13").Value = 46070

# Row 14
$ws.Range("C14").Value = 46070

# Row 15
$ws.Range("C15").Value = 46070

# Row 16
$ws.Range("C16").Value = 46070

# Row 17
$ws.Range("C17").Value = 46070

# Row 18
$ws.Range("A18").Value = 'A 12977-2025'
$ws.Range("B18").Value = 45734.45465277778
$ws.Range("C18").Value = 46070
$ws.Range("G18").Value = 2.1

# Row 19
$ws.Range("A19").Value = 'A 61167-2024'
$ws.Range("B19").Value = 45645
$ws.Range("C19").Value = 46070
$ws.Range("G19").Value = 3

# Row 20
$ws.Range("C20").Value = 46070

# Row 21
$ws.Range("A21").Value = 'A 43067-2024'
$ws.Range("B21").Value = 45567.47446759259
$ws.Range("C21").Value = 46070
$ws.Range("G21").Value = 1.1

# Row 22
$ws.Range("A22").Value = 'A 62768-2025'
$ws.Range("B22").Value = 46008.59856481481
$ws.Range("C22").Value = 46070
$ws.Range("G22").Value = 4.2

# Row 23
$ws.Range("C23").Value = 46070

# Row 24
$ws.Range("A24").Value = 'A 44926-2025'
$ws.Range("B24").Value = 45918.55856481481
$ws.Range("C24").Value = 46070
$ws.Range("G24").Value = 3.8

# Row 25
$ws.Range("A25").Value = 'A 33201-2023'
$ws.Range("B25").Value = 45127.42379629629
$ws.Range("C25").Value = 46070
$ws.Range("G25").Value = 0.9

# Row 26
$ws.Range("A26").Value = 'A 6042-2024'
$ws.Range("B26").Value = 45336
$ws.Range("C26").Value = 46070
$ws.Range("G26").Value = 1.7

# Row 27
$ws.Range("A27").Value = 'A 37570-2025'
$ws.Range("B27").Value = 45880.37358796296
$ws.Range("C27").Value = 46070
$ws.Range("G27").Value = 0.9

# Row 28
$ws.Range("A28").Value = 'A 61178-2024'
$ws.Range("B28").Value = 45645
$ws.Range("C28").Value = 46070
$ws.Range("G28").Value = 6.6

# Row 29
$ws.Range("A29").Value = 'A 46384-2025'
$ws.Range("B29").Value = 45925
$ws.Range("C29").Value = 46070
$ws.Range("G29").Value = 3.1

# Row 30
$ws.Range("A30").Value = 'A 31321-2025'
$ws.Range("B30").Value = 45833
$ws.Range("C30").Value = 46070
$ws.Range("G30").Value = 6.4

# Row 31
$ws.Range("A31").Value = 'A 48265-2025'
$ws.Range("B31").Value = 45933
$ws.Range("C31").Value = 46070
$ws.Range("G31").Value = 2.1

# Row 32
$ws.Range("A32").Value = 'A 22953-2023'
$ws.Range("B32").Value = 45072
$ws.Range("C32").Value = 46070
$ws.Range("G32").Value = 1.9

# Row 33
$ws.Range("A33").Value = 'A 46379-2025'
$ws.Range("B33").Value = 45925
$ws.Range("C33").Value = 46070
$ws.Range("G33").Value = 7.1

# Row 34
$ws.Range("A34").Value = 'A 62804-2023'
$ws.Range("B34").Value = 45271
$ws.Range("C34").Value = 46070
$ws.Range("G34").Value = 0.6

# Row 35
$ws.Range("A35").Value = 'A 3811-2024'
$ws.Range("B35").Value = 45321.673125
$ws.Range("C35").Value = 46070
$ws.Range("G35").Value = 0.9

# Row 36
$ws.Range("A36").Value = 'A 3676-2022'
$ws.Range("B36").Value = 44586
$ws.Range("C36").Value = 46070
$ws.Range("G36").Value = 0.5

# Row 37
$ws.Range("A37").Value = 'A 6258-2024'
$ws.Range("B37").Value = 45337.77947916667
$ws.Range("C37").Value = 46070
$ws.Range("G37").Value = 1.4

# Row 38
$ws.Range("A38").Value = 'A 88-2025'
$ws.Range("B38").Value = 45659.46386574074
$ws.Range("C38").Value = 46070
$ws.Range("G38").Value = 1.4

# Row 39
$ws.Range("A39").Value = 'A 8639-2023'
$ws.Range("B39").Value = 44977.95614583333
$ws.Range("C39").Value = 46070
$ws.Range("G39").Value = 3.8

# Row 40
$ws.Range("A40").Value = 'A 60809-2024'
$ws.Range("B40").Value = 45644.61414351852
$ws.Range("C40").Value = 46070
$ws.Range("G40").Value = 0.5

# Row 41
$ws.Range("C41").Value = 46070

# Row 42
$ws.Range("C42").Value = 46070

# Row 43
$ws.Range("C43").Value = 46070

# Row 44
$ws.Range("C44").Value = 46070

# Row 45
$ws.Range("C45").Value = 46070

# Row 46
$ws.Range("C46").Value = 46070

# Row 47
$ws.Range("C47").Value = 46070

